$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text, matching the source data which uses
# localized thousand-separator formatted strings (e.g. "91.114.59").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "91.114.59"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "3.172.80"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "215.97"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "629.35"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").Value = "1.15"
$ws.Range("E7").Value = "  +30.76%  "
$ws.Range("D8").Value = "0.373"
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "3.169.87"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").Value = "0.761"
$ws.Range("E11").Value = "  +13.63%  "
$ws.Range("E12").Value = "  +8.01%  "
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("D14").Value = "5.68"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").Value = "35.08"
$ws.Range("E15").Value = "  +9.11%  "
$ws.Range("D16").Value = "90.603.85"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "3.754.60"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").Value = "3.130.52"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("E19").Value = "  +11.94%  "
$ws.Range("D20").Value = "14.49"
$ws.Range("E20").Value = "  +8.59%  "
$ws.Range("D21").Value = "468.39"
$ws.Range("E21").Value = "  +10.86%  "
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "9.14"
$ws.Range("E23").Value = "  +11.03%  "
$ws.Range("D24").Value = "5.24"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").Value = "5.93"
$ws.Range("E25").Value = "  +10.92%  "
$ws.Range("D26").Value = "94.80"
$ws.Range("E26").Value = "  +14.10%  "
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").Value = "  +6.11%  "
$ws.Range("D28").Value = "3.331.65"
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "9.32"
$ws.Range("E30").Value = "  +12.66%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "0.163"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "27.69"
$ws.Range("E33").Value = "  +22.22%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.194"
$ws.Range("E34").Value = "  +43.46%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "526.75"
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("E36").Value = "  +8.02%  "
$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  +5.14%  "
$ws.Range("D39").Value = "0.144"
$ws.Range("E39").Value = "  +8.97%  "
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("D41").Value = "22.25"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "0.0863"
$ws.Range("E42").Value = "  +26.36%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "0.415"
$ws.Range("E44").Value = "  +15.48%  "
$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  +9.28%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "151.53"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "4.62"
$ws.Range("E48").Value = "  +10.24%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.688"
$ws.Range("E49").Value = "  +17.28%  "
$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  +12.66%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "45.34"
$ws.Range("E51").Value = "  +4.25%  "
